$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.700.79"
$ws.Range("E2").Value = "  +3.59%  "

$ws.Range("D3").Value = "2.347.55"
$ws.Range("E3").Value = "  +3.02%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.36"
$ws.Range("E5").Value = "  +2.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.48"
$ws.Range("E6").Value = "  +3.75%  "

$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  +1.10%  "

$ws.Range("D9").Value = "2.343.52"
$ws.Range("E9").Value = "  +2.27%  "

$ws.Range("E10").Value = "  +6.56%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +6.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.64"
$ws.Range("E14").Value = "  +1.33%  "

$ws.Range("D15").Value = "2.737.92"
$ws.Range("E15").Value = "  +2.06%  "

$ws.Range("D16").Value = "56.702.97"
$ws.Range("E16").Value = "  +3.58%  "

$ws.Range("E17").Value = "  +2.34%  "

$ws.Range("D18").Value = "2.346.64"
$ws.Range("E18").Value = "  +2.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.38"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("E20").Value = "  +2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "319.15"
$ws.Range("E21").Value = "  +4.01%  "

$ws.Range("E22").Value = "  +3.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.63"
$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.159"
$ws.Range("E26").Value = "  +5.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.72"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.92"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.22"
$ws.Range("E29").Value = "  +10.06%  "

$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +4.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.22"
$ws.Range("E31").Value = "  +2.56%  "

$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.24"
$ws.Range("E33").Value = "  +1.69%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("E36").Value = "  +3.32%  "

$ws.Range("E37").Value = "  +3.15%  "

$ws.Range("E38").Value = "  +4.92%  "

$ws.Range("E39").Value = "  +7.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.39"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.85"
$ws.Range("E42").Value = "  +9.46%  "

$ws.Range("E43").Value = "  +4.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.38"
$ws.Range("E44").Value = "  +10.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.03"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("E46").Value = "  +2.93%  "

$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.559"
$ws.Range("E48").Value = "  +1.71%  "

$ws.Range("E49").Value = "  +4.24%  "

$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("E51").Value = "  +1.81%  "
